$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6616.0415
$ws.Range("I51").Value = 2245.25
$ws.Range("J51").Value = 7490.2
$ws.Range("K51").Value = 2245.25
$ws.Range("L51").Value = 7490.2
$ws.Range("M51").Value = -1761.25
$ws.Range("N51").Value = -8458.200000000001
$ws.Range("H62").Value = 9939.9
$ws.Range("I62").Value = 800
$ws.Range("K62").Value = 800
$ws.Range("M62").Value = -176
$ws.Range("H65").Value = 9939.9
$ws.Range("I65").Value = 800
$ws.Range("K65").Value = 4000
$ws.Range("M65").Value = -880
$ws.Range("H76").Value = 5564.6523
$ws.Range("I76").Value = 4578.846
$ws.Range("J76").Value = 6846.2
$ws.Range("K76").Value = 4578.846
$ws.Range("L76").Value = 6846.2
$ws.Range("M76").Value = -4263.846
$ws.Range("N76").Value = -7476.2
$ws.Range("H79").Value = 5564.6523
$ws.Range("I79").Value = 4578.846
$ws.Range("J79").Value = 6846.2
$ws.Range("K79").Value = 4578.846
$ws.Range("L79").Value = 6846.2
$ws.Range("M79").Value = -3486.846
$ws.Range("N79").Value = -9030.200000000001
$ws.Range("H100").Value = 2779.2727
$ws.Range("I100").Value = 2779.2727
$ws.Range("K100").Value = 2779.2727
$ws.Range("M100").Value = -2238.2727
$ws.Range("H107").Value = 1912.6666
$ws.Range("I107").Value = 1912.6666
$ws.Range("K107").Value = 1912.6666
$ws.Range("M107").Value = 7.333399999999983
$ws.Range("H112").Value = 7630.5884
$ws.Range("J112").Value = 7630.5884
$ws.Range("L112").Value = 22891.7652
$ws.Range("N112").Value = -25107.7652
$ws.Range("H138").Value = 1847.56
$ws.Range("I138").Value = 936.6667
$ws.Range("J138").Value = 2688.3845
$ws.Range("K138").Value = 2810.0001
$ws.Range("L138").Value = 8065.1535
$ws.Range("M138").Value = 2329.9999
$ws.Range("N138").Value = -18345.1535
$ws.Range("H141").Value = 1370.973
$ws.Range("I141").Value = 1096
$ws.Range("K141").Value = 3288
$ws.Range("M141").Value = 1892
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6120.415
$ws.Range("I32").Value = 4575.9707
$ws.Range("J32").Value = 8884.157999999999
$ws.Range("K32").Value = 4575.9707
$ws.Range("L32").Value = 8884.157999999999
$ws.Range("M32").Value = -4288.9707
$ws.Range("N32").Value = -9458.157999999999
$ws.Range("H74").Value = 147331
$ws.Range("I74").Value = 112402.555
$ws.Range("J74").Value = 225920
$ws.Range("K74").Value = 112402.555
$ws.Range("L74").Value = 225920
$ws.Range("M74").Value = -111528.555
$ws.Range("N74").Value = -227668
$ws.Range("H77").Value = 147331
$ws.Range("I77").Value = 112402.555
$ws.Range("J77").Value = 225920
$ws.Range("K77").Value = 562012.7749999999
$ws.Range("L77").Value = 1129600
$ws.Range("M77").Value = -557644.7749999999
$ws.Range("N77").Value = -1138336
$ws.Range("H132").Value = 1512.9111
$ws.Range("I132").Value = 1124.0731
$ws.Range("K132").Value = 3372.2193
$ws.Range("M132").Value = -842.2193000000002
$ws.Range("H139").Value = 98430.5
$ws.Range("J139").Value = 98430.5
$ws.Range("L139").Value = 98430.5
$ws.Range("N139").Value = -108710.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1639.32
$ws.Range("I105").Value = 1359.9
$ws.Range("K105").Value = 1359.9
$ws.Range("M105").Value = 387.0999999999999
$ws.Range("H134").Value = 4375.273
$ws.Range("I134").Value = 1230.9286
$ws.Range("K134").Value = 3692.7858
$ws.Range("M134").Value = -1157.7858
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4000
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4224
$ws.Range("H31").Value = 2993.932
$ws.Range("I31").Value = 2480.1365
$ws.Range("J31").Value = 3507.7273
$ws.Range("K31").Value = 2480.1365
$ws.Range("L31").Value = 3507.7273
$ws.Range("M31").Value = -2185.1365
$ws.Range("N31").Value = -4097.7273
$ws.Range("H34").Value = 2993.932
$ws.Range("I34").Value = 2480.1365
$ws.Range("J34").Value = 3507.7273
$ws.Range("K34").Value = 2480.1365
$ws.Range("L34").Value = 3507.7273
$ws.Range("M34").Value = -2278.1365
$ws.Range("N34").Value = -3911.7273
$ws.Range("H132").Value = 1742.1666
$ws.Range("I132").Value = 1836.909
$ws.Range("J132").Value = 700
$ws.Range("K132").Value = 5510.727000000001
$ws.Range("L132").Value = 2100
$ws.Range("M132").Value = -2980.727000000001
$ws.Range("N132").Value = -7160
$ws.Range("H134").Value = 31540.188
$ws.Range("I134").Value = 41785.434
$ws.Range("J134").Value = 5357.8887
$ws.Range("K134").Value = 125356.302
$ws.Range("L134").Value = 16073.6661
$ws.Range("M134").Value = -122821.302
$ws.Range("N134").Value = -21143.6661
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 6250342.5
$ws.Range("I2").Value = 74.59999999999999
$ws.Range("J2").Value = 7407800
$ws.Range("K2").Value = 447.6
$ws.Range("L2").Value = 44446800
$ws.Range("M2").Value = -334.6
$ws.Range("N2").Value = -44447026
$ws.Range("H37").Value = 43700
$ws.Range("J37").Value = 43700
$ws.Range("L37").Value = 131100
$ws.Range("N37").Value = -131324
$ws.Range("H56").Value = 9265579
$ws.Range("I56").Value = 9265579
$ws.Range("K56").Value = 9265579
$ws.Range("M56").Value = -9265049
$ws.Range("H113").Value = 3020.2812
$ws.Range("I113").Value = 4471.5835
$ws.Range("K113").Value = 13414.7505
$ws.Range("M113").Value = -11244.7505
$ws.Range("H117").Value = 1376.3636
$ws.Range("J117").Value = 1407.75
$ws.Range("L117").Value = 4223.25
$ws.Range("N117").Value = -11107.25
$ws.Range("H121").Value = 11111758
$ws.Range("I121").Value = 409.75
$ws.Range("K121").Value = 1229.25
$ws.Range("M121").Value = 80.75
$ws.Range("H139").Value = 2308.25
$ws.Range("I139").Value = 924.75
$ws.Range("K139").Value = 2774.25
$ws.Range("M139").Value = 2365.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2999
$ws.Range("I5").Value = 998
$ws.Range("K5").Value = 998
$ws.Range("M5").Value = -886
$ws.Range("H122").Value = 172825.73
$ws.Range("I122").Value = 234510.38
$ws.Range("K122").Value = 703531.14
$ws.Range("M122").Value = -701081.14
$ws.Range("H126").Value = 3132.682
$ws.Range("I126").Value = 2830.2354
$ws.Range("J126").Value = 4161
$ws.Range("K126").Value = 8490.706200000001
$ws.Range("L126").Value = 12483
$ws.Range("M126").Value = -6020.706200000001
$ws.Range("N126").Value = -17423
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2509250
$ws.Range("J2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("N2").Value = -15224
$ws.Range("H40").Value = 12245.875
$ws.Range("I40").Value = 9656.333000000001
$ws.Range("K40").Value = 9656.333000000001
$ws.Range("M40").Value = -9520.333000000001
$ws.Range("H93").Value = 1904.5518
$ws.Range("I93").Value = 1992.7826
$ws.Range("K93").Value = 1992.7826
$ws.Range("M93").Value = -744.7826
$ws.Range("H132").Value = 6399.392
$ws.Range("J132").Value = 7817.8184
$ws.Range("L132").Value = 23453.4552
$ws.Range("N132").Value = -28513.4552
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -9888
$ws.Range("N2").ClearContents()
$ws.Range("H41").Value = 77474.92999999999
$ws.Range("J41").Value = 77474.92999999999
$ws.Range("L41").Value = 77474.92999999999
$ws.Range("N41").Value = -78254.92999999999
$ws.Range("H96").Value = 1301
$ws.Range("J96").Value = 1654.75
$ws.Range("L96").Value = 1654.75
$ws.Range("N96").Value = -4400.75
$ws.Range("H132").Value = 9184.299999999999
$ws.Range("I132").Value = 9482.611000000001
$ws.Range("J132").Value = 6499.5
$ws.Range("K132").Value = 28447.833
$ws.Range("L132").Value = 19498.5
$ws.Range("M132").Value = -25917.833
$ws.Range("N132").Value = -24558.5
$ws.Range("H136").Value = 2059.5957
$ws.Range("I136").Value = 1372.7941
$ws.Range("J136").Value = 3855.8462
$ws.Range("K136").Value = 4118.3823
$ws.Range("L136").Value = 11567.5386
$ws.Range("M136").Value = -1568.3823
$ws.Range("N136").Value = -16667.5386
